$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 15.34689499495649
$ws.Range("C2").Value = 11.38288564936875
$ws.Range("D2").Value = 9.919769039986615
$ws.Range("F2").Value = 31.36549369557103
$ws.Range("G2").Value = 31.55481353683188
$ws.Range("H2").Value = 15.01728774433545
$ws.Range("J2").Value = 10.76686407981158
$ws.Range("L2").Value = 11.46242400502128
$ws.Range("M2").Value = 15.9536707943452
$ws.Range("N2").Value = 18.53687835880071
$ws.Range("O2").Value = 23.22858824764555
# Row 3
$ws.Range("B3").Value = 14.95299315747572
$ws.Range("C3").Value = 11.29006046405547
$ws.Range("D3").Value = 9.929853563220792
$ws.Range("F3").Value = 31.43219209273901
$ws.Range("G3").Value = 31.59076113083892
$ws.Range("H3").Value = 15.06007942729439
$ws.Range("J3").Value = 10.78885455767806
$ws.Range("L3").Value = 11.46530572362861
$ws.Range("M3").Value = 15.87389583496579
$ws.Range("N3").Value = 18.58329588922966
$ws.Range("O3").Value = 23.29109539413584
# Row 4
$ws.Range("B4").Value = 14.70759364190124
$ws.Range("C4").Value = 11.23255469776463
$ws.Range("D4").Value = 9.937202286959449
$ws.Range("F4").Value = 31.48007770664595
$ws.Range("G4").Value = 31.62230380546203
$ws.Range("H4").Value = 15.08870589780151
$ws.Range("J4").Value = 10.8030563511736
$ws.Range("L4").Value = 11.46826813968228
$ws.Range("M4").Value = 15.82658155879067
$ws.Range("N4").Value = 18.61353545146029
$ws.Range("O4").Value = 23.3342955274351
# Row 5
$ws.Range("B5").Value = 14.60684485521293
$ws.Range("C5").Value = 11.20900528955855
$ws.Range("D5").Value = 9.940488376702969
$ws.Range("F5").Value = 31.50133208115517
$ws.Range("G5").Value = 31.63753414213387
$ws.Range("H5").Value = 15.10096283324634
$ws.Range("J5").Value = 10.80902010204168
$ws.Range("L5").Value = 11.46977631940609
$ws.Range("M5").Value = 15.80773354722901
$ws.Range("N5").Value = 18.62629647113368
$ws.Range("O5").Value = 23.35311017306278
# Row 6
$ws.Range("B6").Value = 14.59007487268287
$ws.Range("C6").Value = 11.20508834428785
$ws.Range("D6").Value = 9.941051645936149
$ws.Range("F6").Value = 31.50496639690443
$ws.Range("G6").Value = 31.64020648562424
$ws.Range("H6").Value = 15.10303380601407
$ws.Range("J6").Value = 10.81002104736634
$ws.Range("L6").Value = 11.47004495679645
$ws.Range("M6").Value = 15.8046303832509
$ws.Range("N6").Value = 18.62844191769997
$ws.Range("O6").Value = 23.35630736737591
# Row 7
$ws.Range("B7").Value = 14.70623773491159
$ws.Range("C7").Value = 11.23223755193976
$ws.Range("D7").Value = 9.937245423647594
$ws.Range("F7").Value = 31.48035730673907
$ws.Range("G7").Value = 31.62249959265256
$ws.Range("H7").Value = 15.08886880459472
$ws.Range("J7").Value = 10.80313606548549
$ws.Range("L7").Value = 11.46828725965184
$ws.Range("M7").Value = 15.82632559745973
$ws.Range("N7").Value = 18.61370577581673
$ws.Range("O7").Value = 23.33454437075061
# Row 8
$ws.Range("B8").Value = 15.21189694672295
$ws.Range("C8").Value = 11.35099336676615
$ws.Range("D8").Value = 9.923006427168591
$ws.Range("F8").Value = 31.38705097919825
$ws.Range("G8").Value = 31.56523998367903
$ws.Range("H8").Value = 15.0315540826145
$ws.Range("J8").Value = 10.77430150468183
$ws.Range("L8").Value = 11.46317060743357
$ws.Range("M8").Value = 15.92582649429184
$ws.Range("N8").Value = 18.55252270977815
$ws.Range("O8").Value = 23.24913911311352
# Row 9
$ws.Range("B9").Value = 16.16945738961407
$ws.Range("C9").Value = 11.57925032548489
$ws.Range("D9").Value = 9.904237819485367
$ws.Range("F9").Value = 31.25918793201898
$ws.Range("G9").Value = 31.52826692893303
$ws.Range("H9").Value = 14.93782636347904
$ws.Range("J9").Value = 10.72328481128765
$ws.Range("L9").Value = 11.46255899591854
$ws.Range("M9").Value = 16.13356567781354
$ws.Range("N9").Value = 18.44630229293928
$ws.Range("O9").Value = 23.11998837694161
# Row 10
$ws.Range("B10").Value = 16.84489489701926
$ws.Range("C10").Value = 11.7433767492065
$ws.Range("D10").Value = 9.895994930359512
$ws.Range("F10").Value = 31.19896414420071
$ws.Range("G10").Value = 31.54715688940209
$ws.Range("H10").Value = 14.8803471170599
$ws.Range("J10").Value = 10.68914055296026
$ws.Range("L10").Value = 11.46779267895881
$ws.Range("M10").Value = 16.29303209399178
$ws.Range("N10").Value = 18.37659660992945
$ws.Range("O10").Value = 23.0485693534891
# Row 11
$ws.Range("B11").Value = 17.1447004114751
$ws.Range("C11").Value = 11.81710778690536
$ws.Range("D11").Value = 9.893441372831779
$ws.Range("F11").Value = 31.17890504191723
$ws.Range("G11").Value = 31.56575011901681
$ws.Range("H11").Value = 14.85667033760897
$ws.Range("J11").Value = 10.67432534770487
$ws.Range("L11").Value = 11.47139349137345
$ws.Range("M11").Value = 16.36687349301654
$ws.Range("N11").Value = 18.34668429858097
$ws.Range("O11").Value = 23.02119317898249
# Row 12
$ws.Range("B12").Value = 17.2570550998593
$ws.Range("C12").Value = 11.84488077517164
$ws.Range("D12").Value = 9.892645662164345
$ws.Range("F12").Value = 31.17236477870606
$ws.Range("G12").Value = 31.57422620802778
$ws.Range("H12").Value = 14.84805989123994
$ws.Range("J12").Value = 10.66881783424334
$ws.Range("L12").Value = 11.47293113400119
$ws.Range("M12").Value = 16.39500566107522
$ws.Range("N12").Value = 18.33561490177897
$ws.Range("O12").Value = 23.01156291595801
# Row 13
$ws.Range("B13").Value = 17.23291132057576
$ws.Range("C13").Value = 11.83890613141426
$ws.Range("D13").Value = 9.892809427867773
$ws.Range("F13").Value = 31.1737263848848
$ws.Range("G13").Value = 31.57233695897906
$ws.Range("H13").Value = 14.84989849487397
$ws.Range("J13").Value = 10.66999941603065
$ws.Range("L13").Value = 11.4725922543233
$ws.Range("M13").Value = 16.38893960318486
$ws.Range("N13").Value = 18.33798744177265
$ws.Range("O13").Value = 23.01360419088421
# Row 14
$ws.Range("B14").Value = 17.15396801036245
$ws.Range("C14").Value = 11.81939571710009
$ws.Range("D14").Value = 9.893372481212397
$ws.Range("F14").Value = 31.17834581122062
$ws.Range("G14").Value = 31.56641870915863
$ws.Range("H14").Value = 14.8559548268325
$ws.Range("J14").Value = 10.67387018588238
$ws.Range("L14").Value = 11.47151651297258
$ws.Range("M14").Value = 16.36918463281246
$ws.Range("N14").Value = 18.3457684508676
$ws.Range("O14").Value = 23.02038612418109
# Row 15
$ws.Range("B15").Value = 17.1054570041731
$ws.Range("C15").Value = 11.80742543818319
$ws.Range("D15").Value = 9.893739648786386
$ws.Range("F15").Value = 31.18131283458814
$ws.Range("G15").Value = 31.56298041129176
$ws.Range("H15").Value = 14.85971079697433
$ws.Range("J15").Value = 10.67625450335493
$ws.Range("L15").Value = 11.47088022155038
$ws.Range("M15").Value = 16.35710578344324
$ws.Range("N15").Value = 18.35056808984996
$ws.Range("O15").Value = 23.0246362032459
# Row 16
$ws.Range("B16").Value = 16.82514287007629
$ws.Range("C16").Value = 11.73853836213428
$ws.Range("D16").Value = 9.896185818970283
$ws.Range("F16").Value = 31.20042275832206
$ws.Range("G16").Value = 31.54614280541368
$ws.Range("H16").Value = 14.88194418151523
$ws.Range("J16").Value = 10.69012315478123
$ws.Range("L16").Value = 11.4675817936047
$ws.Range("M16").Value = 16.28823109941297
$ws.Range("N16").Value = 18.37858755442382
$ws.Range("O16").Value = 23.05046144731925
# Row 17
$ws.Range("B17").Value = 16.65119365971677
$ws.Range("C17").Value = 11.69603028024259
$ws.Range("D17").Value = 9.89799231962524
$ws.Range("F17").Value = 31.21402586698281
$ws.Range("G17").Value = 31.53837327901303
$ws.Range("H17").Value = 14.89621658155131
$ws.Range("J17").Value = 10.69881449271946
$ws.Range("L17").Value = 11.46586984983428
$ws.Range("M17").Value = 16.24629962084392
$ws.Range("N17").Value = 18.39623639145717
$ws.Range("O17").Value = 23.06761487337022
# Row 18
$ws.Range("B18").Value = 16.55044662124611
$ws.Range("C18").Value = 11.67149393029594
$ws.Range("D18").Value = 9.899144007728635
$ws.Range("F18").Value = 31.22254060203341
$ws.Range("G18").Value = 31.53484597465737
$ws.Range("H18").Value = 14.90465824680293
$ws.Range("J18").Value = 10.70388105216363
$ws.Range("L18").Value = 11.46500009991909
$ws.Range("M18").Value = 16.22230506530972
$ws.Range("N18").Value = 18.4065567319507
$ws.Range("O18").Value = 23.0779622367333
# Row 19
$ws.Range("B19").Value = 16.51621915072492
$ws.Range("C19").Value = 11.66317184603545
$ws.Range("D19").Value = 9.899553316896901
$ws.Range("F19").Value = 31.22554212655911
$ws.Range("G19").Value = 31.53381345556879
$ws.Range("H19").Value = 14.90755638725165
$ws.Range("J19").Value = 10.70560811448612
$ws.Range("L19").Value = 11.46472539401829
$ws.Range("M19").Value = 16.21420260983021
$ws.Range("N19").Value = 18.41008010261344
$ws.Range("O19").Value = 23.08154827396974
# Row 20
$ws.Range("B20").Value = 16.6697836653784
$ws.Range("C20").Value = 11.70056442267068
$ws.Range("D20").Value = 9.897788362315744
$ws.Range("F20").Value = 31.21250631354155
$ws.Range("G20").Value = 31.53910292937841
$ws.Range("H20").Value = 14.89467319024439
$ws.Range("J20").Value = 10.69788229893336
$ws.Range("L20").Value = 11.46604020541114
$ws.Range("M20").Value = 16.25075065595059
$ws.Range("N20").Value = 18.39434013644639
$ws.Range("O20").Value = 23.06573905288234
# Row 21
$ws.Range("B21").Value = 17.17718824548341
$ws.Range("C21").Value = 11.82513050436081
$ws.Range("D21").Value = 9.893202456812569
$ws.Range("F21").Value = 31.17696032023374
$ws.Range("G21").Value = 31.56811812058322
$ws.Range("H21").Value = 14.85416628791751
$ws.Range("J21").Value = 10.67273046370212
$ws.Range("L21").Value = 11.47182777043924
$ws.Range("M21").Value = 16.37498266695225
$ws.Range("N21").Value = 18.34347599006852
$ws.Range("O21").Value = 23.01837410824417
# Row 22
$ws.Range("B22").Value = 17.50191519306015
$ws.Range("C22").Value = 11.9056766409065
$ws.Range("D22").Value = 9.891203197662113
$ws.Range("F22").Value = 31.15988222677138
$ws.Range("G22").Value = 31.59544447818694
$ws.Range("H22").Value = 14.82976448734768
$ws.Range("J22").Value = 10.6568905970508
$ws.Range("L22").Value = 11.47662435979326
$ws.Range("M22").Value = 16.45715859916504
$ws.Range("N22").Value = 18.31173536899925
$ws.Range("O22").Value = 22.99171171965271
# Row 23
$ws.Range("B23").Value = 17.32926464930335
$ws.Range("C23").Value = 11.86277119912692
$ws.Range("D23").Value = 9.892179188825938
$ws.Range("F23").Value = 31.16843402359805
$ws.Range("G23").Value = 31.58009597221497
$ws.Range("H23").Value = 14.84259857827001
$ws.Range("J23").Value = 10.66529003328343
$ws.Range("L23").Value = 11.47397200039941
$ws.Range("M23").Value = 16.41321541173782
$ws.Range("N23").Value = 18.32853871991398
$ws.Range("O23").Value = 23.00554872188939
# Row 24
$ws.Range("B24").Value = 16.66138142491657
$ws.Range("C24").Value = 11.69851484076092
$ws.Range("D24").Value = 9.897880219054713
$ws.Range("F24").Value = 31.21319114155811
$ws.Range("G24").Value = 31.53877012779256
$ws.Range("H24").Value = 14.89537022148768
$ws.Range("J24").Value = 10.69830352630965
$ws.Range("L24").Value = 11.4659628309983
$ws.Range("M24").Value = 16.24873799155262
$ws.Range("N24").Value = 18.39519689194049
$ws.Range("O24").Value = 23.06658559871885
# Row 25
$ws.Range("B25").Value = 15.91484765764048
$ws.Range("C25").Value = 11.51807405927755
$ws.Range("D25").Value = 9.908338589944067
$ws.Range("F25").Value = 31.28786414253409
$ws.Range("G25").Value = 31.53018592418818
$ws.Range("H25").Value = 14.96118346422204
$ws.Range("J25").Value = 10.73649773234713
$ws.Range("L25").Value = 11.46172179433042
$ws.Range("M25").Value = 16.07610215995857
$ws.Range("N25").Value = 18.47357030472904
$ws.Range("O25").Value = 23.15081250757351
